$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = "system, backup@backdoor.com, System"
$v = $ws.Range("G2").Value2
Write-Host "VAL:" $v
